$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -18.11838856361434
$ws.Range("C2").Value = 1.905415864348617
$ws.Range("D2").Value = -18.11838856361434
$ws.Range("E2").Value = -18.11838856361434
$ws.Range("F2").Value = -18.11838856361434
$ws.Range("G2").Value = -18.11838856361434
$ws.Range("H2").Value = -18.11838856361434
$ws.Range("I2").Value = -18.11838856361434
$ws.Range("J2").Value = -18.11838856361434
$ws.Range("K2").Value = -18.11838856361434
$ws.Range("B3").Value = -18.11838856361434
$ws.Range("C3").Value = -18.11838856361434
$ws.Range("D3").Value = -18.11838856361434
$ws.Range("E3").Value = -18.11838856361434
$ws.Range("F3").Value = -18.11838856361434
$ws.Range("G3").Value = -18.11838856361434
$ws.Range("H3").Value = -18.11838856361434
$ws.Range("I3").Value = 4.321923278510755
$ws.Range("J3").Value = -18.11838856361434
$ws.Range("K3").Value = -18.11838856361434
$ws.Range("B4").Value = -18.11838856361434
$ws.Range("C4").Value = 2.190738025380502
$ws.Range("D4").Value = 2.155849661594278
$ws.Range("E4").Value = -18.11838856361434
$ws.Range("F4").Value = 3.438338722808441
$ws.Range("G4").Value = -18.11838856361434
$ws.Range("H4").Value = 1.213763838496409
$ws.Range("I4").Value = -18.11838856361434
$ws.Range("J4").Value = 2.330305499617314
$ws.Range("K4").Value = -18.11838856361434
$ws.Range("B5").Value = -18.11838856361434
$ws.Range("C5").Value = 1.865207697283254
$ws.Range("D5").Value = -18.11838856361434
$ws.Range("E5").Value = -18.11838856361434
$ws.Range("F5").Value = -18.11838856361434
$ws.Range("G5").Value = 2.992971924230173
$ws.Range("H5").Value = -18.11838856361434
$ws.Range("I5").Value = -18.11838856361434
$ws.Range("J5").Value = -18.11838856361434
$ws.Range("K5").Value = -18.11838856361434
$ws.Range("B6").Value = -18.11838856361434
$ws.Range("C6").Value = -18.11838856361434
$ws.Range("D6").Value = -18.11838856361434
$ws.Range("E6").Value = -18.11838856361434
$ws.Range("F6").Value = -18.11838856361434
$ws.Range("G6").Value = -18.11838856361434
$ws.Range("H6").Value = -18.11838856361434
$ws.Range("I6").Value = -18.11838856361434
$ws.Range("J6").Value = -18.11838856361434
$ws.Range("K6").Value = -18.11838856361434
$ws.Range("B7").Value = 2.669521167669268
$ws.Range("C7").Value = -18.11838856361434
$ws.Range("D7").Value = -18.11838856361434
$ws.Range("E7").Value = -18.11838856361434
$ws.Range("F7").Value = -18.11838856361434
$ws.Range("G7").Value = -18.11838856361434
$ws.Range("H7").Value = -18.11838856361434
$ws.Range("I7").Value = -18.11838856361434
$ws.Range("J7").Value = -18.11838856361434
$ws.Range("K7").Value = -18.11838856361434
$ws.Range("B8").Value = -18.11838856361434
$ws.Range("C8").Value = -18.11838856361434
$ws.Range("D8").Value = -18.11838856361434
$ws.Range("E8").Value = 1.754648447770726
$ws.Range("F8").Value = -18.11838856361434
$ws.Range("G8").Value = -18.11838856361434
$ws.Range("H8").Value = -18.11838856361434
$ws.Range("I8").Value = -18.11838856361434
$ws.Range("J8").Value = -18.11838856361434
$ws.Range("K8").Value = -18.11838856361434
$ws.Range("B9").Value = 3.769534462107125
$ws.Range("C9").Value = -18.11838856361434
$ws.Range("D9").Value = -18.11838856361434
$ws.Range("E9").Value = -18.11838856361434
$ws.Range("F9").Value = -18.11838856361434
$ws.Range("G9").Value = -18.11838856361434
$ws.Range("H9").Value = -18.11838856361434
$ws.Range("I9").Value = -18.11838856361434
$ws.Range("J9").Value = -18.11838856361434
$ws.Range("K9").Value = -18.11838856361434
$ws.Range("B10").Value = -18.11838856361434
$ws.Range("C10").Value = -18.11838856361434
$ws.Range("D10").Value = -18.11838856361434
$ws.Range("E10").Value = -18.11838856361434
$ws.Range("F10").Value = -18.11838856361434
$ws.Range("G10").Value = -18.11838856361434
$ws.Range("H10").Value = -18.11838856361434
$ws.Range("I10").Value = -18.11838856361434
$ws.Range("J10").Value = -18.11838856361434
$ws.Range("K10").Value = 1.766230640683092
$ws.Range("B11").Value = -18.11838856361434
$ws.Range("C11").Value = -18.11838856361434
$ws.Range("D11").Value = -18.11838856361434
$ws.Range("E11").Value = 2.863619710341915
$ws.Range("F11").Value = -18.11838856361434
$ws.Range("G11").Value = 2.521094351183916
$ws.Range("H11").Value = -18.11838856361434
$ws.Range("I11").Value = -18.11838856361434
$ws.Range("J11").Value = -18.11838856361434
$ws.Range("K11").Value = 1.701860529490527
$ws.Range("B12").Value = -18.11838856361434
$ws.Range("C12").Value = -18.11838856361434
$ws.Range("D12").Value = -18.11838856361434
$ws.Range("E12").Value = -18.11838856361434
$ws.Range("F12").Value = -18.11838856361434
$ws.Range("G12").Value = -18.11838856361434
$ws.Range("H12").Value = -18.11838856361434
$ws.Range("I12").Value = -18.11838856361434
$ws.Range("J12").Value = -18.11838856361434
$ws.Range("K12").Value = -18.11838856361434
$ws.Range("B13").Value = -18.11838856361434
$ws.Range("C13").Value = -18.11838856361434
$ws.Range("D13").Value = -18.11838856361434
$ws.Range("E13").Value = 2.415523807158697
$ws.Range("F13").Value = -18.11838856361434
$ws.Range("G13").Value = -18.11838856361434
$ws.Range("H13").Value = -18.11838856361434
$ws.Range("I13").Value = -18.11838856361434
$ws.Range("J13").Value = 2.395505737182692
$ws.Range("K13").Value = 2.096202109026979
$ws.Range("B14").Value = -18.11838856361434
$ws.Range("C14").Value = -18.11838856361434
$ws.Range("D14").Value = 1.241260205594947
$ws.Range("E14").Value = -18.11838856361434
$ws.Range("F14").Value = -18.11838856361434
$ws.Range("G14").Value = -18.11838856361434
$ws.Range("H14").Value = -18.11838856361434
$ws.Range("I14").Value = -18.11838856361434
$ws.Range("J14").Value = -18.11838856361434
$ws.Range("K14").Value = 1.988763367282632
$ws.Range("B15").Value = -18.11838856361434
$ws.Range("C15").Value = -18.11838856361434
$ws.Range("D15").Value = 1.238856289242096
$ws.Range("E15").Value = -18.11838856361434
$ws.Range("F15").Value = -18.11838856361434
$ws.Range("G15").Value = -18.11838856361434
$ws.Range("H15").Value = -18.11838856361434
$ws.Range("I15").Value = -18.11838856361434
$ws.Range("J15").Value = -18.11838856361434
$ws.Range("K15").Value = -18.11838856361434
$ws.Range("B16").Value = -18.11838856361434
$ws.Range("C16").Value = -18.11838856361434
$ws.Range("D16").Value = -18.11838856361434
$ws.Range("E16").Value = -18.11838856361434
$ws.Range("F16").Value = -18.11838856361434
$ws.Range("G16").Value = -18.11838856361434
$ws.Range("H16").Value = -18.11838856361434
$ws.Range("I16").Value = -18.11838856361434
$ws.Range("J16").Value = 2.409429168352498
$ws.Range("K16").Value = -18.11838856361434
$ws.Range("B17").Value = -18.11838856361434
$ws.Range("C17").Value = 2.034946986243954
$ws.Range("D17").Value = 2.286541967420876
$ws.Range("E17").Value = -18.11838856361434
$ws.Range("F17").Value = -18.11838856361434
$ws.Range("G17").Value = -18.11838856361434
$ws.Range("H17").Value = 0.654753863403189
$ws.Range("I17").Value = -18.11838856361434
$ws.Range("J17").Value = 1.19970624698226
$ws.Range("K17").Value = -18.11838856361434
$ws.Range("B18").Value = -18.11838856361434
$ws.Range("C18").Value = -18.11838856361434
$ws.Range("D18").Value = -18.11838856361434
$ws.Range("E18").Value = -18.11838856361434
$ws.Range("F18").Value = -18.11838856361434
$ws.Range("G18").Value = -18.11838856361434
$ws.Range("H18").Value = 1.243574480412999
$ws.Range("I18").Value = -18.11838856361434
$ws.Range("J18").Value = 1.070144746691403
$ws.Range("K18").Value = -18.11838856361434
$ws.Range("B19").Value = -18.11838856361434
$ws.Range("C19").Value = -18.11838856361434
$ws.Range("D19").Value = 1.714799571283544
$ws.Range("E19").Value = -18.11838856361434
$ws.Range("F19").Value = -18.11838856361434
$ws.Range("G19").Value = -18.11838856361434
$ws.Range("H19").Value = 1.652560787140916
$ws.Range("I19").Value = -18.11838856361434
$ws.Range("J19").Value = -18.11838856361434
$ws.Range("K19").Value = -18.11838856361434
$ws.Range("B20").Value = -18.11838856361434
$ws.Range("C20").Value = 0.8818122981385548
$ws.Range("D20").Value = 1.410403181426498
$ws.Range("E20").Value = -18.11838856361434
$ws.Range("F20").Value = 3.195283349125179
$ws.Range("G20").Value = -18.11838856361434
$ws.Range("H20").Value = 2.344909293160006
$ws.Range("I20").Value = -18.11838856361434
$ws.Range("J20").Value = -18.11838856361434
$ws.Range("K20").Value = 2.350586414248165
$ws.Range("B21").Value = -18.11838856361434
$ws.Range("C21").Value = 1.073424910561728
$ws.Range("D21").Value = -18.11838856361434
$ws.Range("E21").Value = 2.004314332221972
$ws.Range("F21").Value = -18.11838856361434
$ws.Range("G21").Value = 2.655041246039333
$ws.Range("H21").Value = 2.463136613190541
$ws.Range("I21").Value = -18.11838856361434
$ws.Range("J21").Value = -18.11838856361434
$ws.Range("K21").Value = -18.11838856361434
